$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Test ID" column (column A, containing the numbers 1/2/3) was
# removed from the demo plan sheet. Deleting the entire column shifts
# App Type / Zip Code / Insurance Type / Plans / Action (formerly
# B:F) left into A:E, carrying their values, styles and column widths
# with them.
$ws.Columns("A").Delete()

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("C8").Select()
